$wb = $excel.ActiveWorkbook

# ----- Summary sheet -----
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1399.81
$wsSummary.Range("B4").Value = -0.2
$wsSummary.Range("B5").Value = -0.21
$wsSummary.Range("B6").Value = 19
$wsSummary.Range("B7").Value = 10
$wsSummary.Range("B9").Value = 52.63

# ----- Strategy Status sheet -----
$wsStrategy = $wb.Worksheets.Item("Strategy Status")
$wsStrategy.Range("C5").Value = 99.81
$wsStrategy.Range("D5").Value = 19
$wsStrategy.Range("E5").Value = -0.2
$wsStrategy.Range("F5").Value = -0.19
$wsStrategy.Range("G5").Value = 52.63

# ----- New trade row data (shared between "All Trades" and "MarketMaking") -----
$newRow = @("2026-02-17", "20:04:19", "MarketMaking", "UP", 0.95, 0.97, "CLOSED", 2.1053, 0.02, 99.81, 0, 0, 0.6, "Normal spread capture: 19600 bps", "early_exit", 0.11)

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Cells.Item(20, 1).Value = 19
    $ws.Cells.Item(20, 2).NumberFormat = "@"
    $ws.Cells.Item(20, 2).Value = $newRow[0]
    $ws.Cells.Item(20, 3).Value = $newRow[1]
    $ws.Cells.Item(20, 4).Value = $newRow[2]
    $ws.Cells.Item(20, 5).Value = $newRow[3]
    $ws.Cells.Item(20, 6).Value = $newRow[4]
    $ws.Cells.Item(20, 7).Value = $newRow[5]
    $ws.Cells.Item(20, 8).Value = $newRow[6]
    $ws.Cells.Item(20, 9).Value = $newRow[7]
    $ws.Cells.Item(20, 10).Value = $newRow[8]
    $ws.Cells.Item(20, 11).Value = $newRow[9]
    $ws.Cells.Item(20, 12).Value = $newRow[10]
    $ws.Cells.Item(20, 13).Value = $newRow[11]
    $ws.Cells.Item(20, 14).Value = $newRow[12]
    $ws.Cells.Item(20, 15).Value = $newRow[13]
    $ws.Cells.Item(20, 16).Value = $newRow[14]
    $ws.Cells.Item(20, 17).Value = $newRow[15]
}
